$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" (column E) list for rows 16-23 gets reversed
# (2011,2012,2101,2102,2103,2104,2105,2106 -> 2106,2105,2104,2103,2102,2101,2012,2011),
# and the "Valor Mora" (column F) values follow the same row-for-row swap.
$periodos = @("2106", "2105", "2104", "2103", "2102", "2101", "2012", "2011")
$valores  = @(24578, 35112, 35112, 35112, 35112, 35112, 35112, 35112)

for ($i = 0; $i -lt 8; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
